# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
# These values mirror the same events' data as refreshed at a later scrape time.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value for "展览" sheet (column F)
$exhibitionUpdates = @{
    5  = 15329
    6  = 415
    8  = 688
    17 = 191
    28 = 32
    29 = 37
    31 = 44
    34 = 294
    37 = 5443
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new value for "全部类型" sheet (column F)
$allTypesUpdates = @{
    5  = 15329
    6  = 415
    8  = 688
    18 = 191
    29 = 32
    30 = 37
    34 = 44
    37 = 294
    40 = 5443
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
